# Add a new entry (row 16) to the coding-questions tracker sheet:
#   - B16: date 12/25/2024 (same date used for the preceding rows)
#   - C16: problem name "longestEvenOddSubarray"
#   - G16: status "solved"
# and leave the selection on the newly filled-in C16 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting already used in column B (row 15) down to B16
# so the new date cell keeps the same number format as the rest of the column.
$null = $ws.Range("B15").Copy()
$null = $ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 45651

# Fill in the new problem name and status; these columns already carry a
# default "wrap text" style (style index 2) via the column definition, so a
# plain value assignment picks it up automatically, same as typing it in.
$ws.Range("C16").Value = "longestEvenOddSubarray"
$ws.Range("G16").Value = "solved"

# Leave the cursor on the cell that was last edited.
$null = $ws.Range("C16").Select()
